$wb = $excel.ActiveWorkbook

# --- market_creation: add row 12 ---
$ws1 = $wb.Worksheets.Item("market_creation")
$ws1.Range("F12").Value = "no"
$ws1.Range("G12").Value = "yes"
$ws1.Range("H12").Value = "no"
$ws1.Range("R12").Value = 0

# --- organisational_change: add row 4 ---
$ws2 = $wb.Worksheets.Item("organisational_change")
$ws2.Range("B4").Value = "yes"
$ws2.Range("C4").Value = "no"
$ws2.Range("D4").Value = "no"
$ws2.Range("Q4").Value = 0

# --- enabling_environment: add row 5 ---
$ws3 = $wb.Worksheets.Item("enabling_environment")
$ws3.Range("B5").Value = "no"
$ws3.Range("C5").Value = "no"
$ws3.Range("D5").Value = "no"
$ws3.Range("O5").Value = 0

# --- user: add rows 40-47 ---
$ws5 = $wb.Worksheets.Item("user")
$newIds = @("3Zpv40", "Wumc41", "XZgj42", "F8lG43", "IETK44", "Np4g45", "oKZK46", "cpFm47")
$r = 40
foreach ($id in $newIds) {
    $ws5.Range("A" + $r).Value = $id
    $ws5.Range("B" + $r).Value = " "
    $ws5.Range("C" + $r).Value = ""
    $ws5.Range("D" + $r).Value = ""
    $ws5.Range("E" + $r).Value = ""
    $ws5.Range("F" + $r).Value = ""
    $ws5.Range("G" + $r).Value = ""
    $ws5.Range("L" + $r).Value = ""
    $r = $r + 1
}
